# Fix character-encoding mojibake in the cover letter.
#
# The document text contains two flavours of mis-decoded UTF-8:
#   "a-circumflex, euro, left-double-quote"  (was originally an EN DASH, U+2013)
#   "a-circumflex, euro, cent"               (was originally a BULLET,   U+2022)
# in both cases the culprit byte sequence starts with U+00E2 ("a-circumflex").
#
# The fix that was actually applied upstream ran a document-wide Find/Replace
# turning every stray U+00E2 into the two characters ">=" -- except inside the
# two "bullet list" paragraphs that were simultaneously restyled from
# "BodyText" to "BlockText" (DECLARATIONS, MANUSCRIPT HIGHLIGHTS), where a
# follow-up paragraph-scoped pass stripped the leading ">" back off again,
# leaving just "=" in front of those bullets.
#
# We rebuild each affected paragraph with Range.InsertXML so the existing
# run layout (one <w:r> per bullet / separator space) is preserved exactly,
# rather than letting a plain Range.Text/Find.Execute edit collapse the
# paragraph down to a single run.

$d = $word.ActiveDocument

$wordMLns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$aCirc = [char]0x00E2   # mis-decoded lead byte ("â")
$euro  = [char]0x20AC   # "€"
$ldq   = [char]0x201C   # "\u201c" (part of the corrupted en dash)
$cent  = [char]0xA2     # "¢"      (part of the corrupted bullet)

$mojiDash   = "$aCirc$euro$ldq"        # "â€“" corrupted EN DASH
$mojiBullet = "$aCirc$euro$cent"       # "â€¢" corrupted BULLET

$fixDash       = ">=$euro$ldq"         # what a stray dash becomes, always
$fixBulletPlus = ">=$euro$cent"        # what a stray bullet becomes (BodyText paragraphs)
$fixBulletEq   = "=$euro$cent"         # what a stray bullet becomes (paragraphs -> BlockText)

function Fix-Runs($runTexts, $bulletReplacement) {
    $out = New-Object System.Collections.ArrayList
    foreach ($t in $runTexts) {
        $u = $t.Replace($mojiDash, $fixDash)
        $u = $u.Replace($mojiBullet, $bulletReplacement)
        [void]$out.Add($u)
    }
    return $out
}

function Build-ParaXml($style, $runTexts) {
    $sb = New-Object System.Text.StringBuilder
    [void]$sb.Append("<w:p $wordMLns><w:pPr><w:pStyle w:val=`"$style`"/></w:pPr>")
    foreach ($t in $runTexts) {
        [void]$sb.Append('<w:r><w:t xml:space="preserve">')
        [void]$sb.Append($t)
        [void]$sb.Append('</w:t></w:r>')
    }
    [void]$sb.Append('</w:p>')
    return $sb.ToString()
}

function Replace-Paragraph($index, $style, $runTexts) {
    $para = $d.Paragraphs.Item($index)
    $xml = Build-ParaXml $style $runTexts
    $para.Range.InsertXML($xml) | Out-Null
}

# --- Paragraph 4: "Re: Submission ... Article â€“" (stays BodyText) ---
$p4 = @(
    "Re: Submission of Original Research Article $mojiDash",
    " ",
    "$ldq",
    "High Consistency, Limited Accuracy: Evaluating Large Language Models for Binary Medical Diagnosis",
    [string][char]0x201D
)
Replace-Paragraph 4 "BodyText" (Fix-Runs $p4 $fixBulletPlus)

# --- Paragraph 11: NOVEL FINDINGS bullets (stays BodyText) ---
$p11 = @(
    "Our study reveals a striking dissociation between consistency and accuracy:",
    " ",
    "$mojiBullet All models achieved exceptional reproducibility (99-100% consistency)",
    " ",
    "$mojiBullet Diagnostic accuracy remained at chance level (~50%)",
    " ",
    "$mojiBullet The consistency-accuracy gap reached ~50 percentage points",
    " ",
    "$mojiBullet Models showed systematic bias toward positive diagnosis (49-51 false positives vs 0-1 false negatives)",
    " ",
    "$mojiBullet Prompt engineering had minimal impact (&lt;3% prediction change)",
    " ",
    "$mojiBullet Error patterns were highly systematic across all three models"
)
Replace-Paragraph 11 "BodyText" (Fix-Runs $p11 $fixBulletPlus)

# --- Paragraph 19: WHY JAMIA bullets (stays BodyText) ---
$p19 = @(
    "This manuscript is an excellent fit for JAMIA because:",
    " ",
    "$mojiBullet Aligns with the journal" + [char]0x2019 + "s focus on AI in medicine and clinical decision support",
    " ",
    "$mojiBullet Addresses timely concerns about LLM reliability in healthcare",
    " ",
    "$mojiBullet Provides rigorous empirical evidence with immediate clinical implications",
    " ",
    "$mojiBullet Appeals to diverse readership (clinicians, AI researchers, policymakers)",
    " ",
    "$mojiBullet Contributes to ongoing dialogue about responsible AI in medicine"
)
Replace-Paragraph 19 "BodyText" (Fix-Runs $p19 $fixBulletPlus)

# --- Paragraph 21: DECLARATIONS bullets -> restyled to BlockText ---
$p21 = @(
    "$mojiBullet This manuscript represents original work not previously published or under consideration elsewhere",
    " ",
    "$mojiBullet A preprint version is available on medRxiv for community feedback and rapid dissemination",
    " ",
    "$mojiBullet All authors have approved the manuscript and agree with submission to JAMIA",
    " ",
    "$mojiBullet We have no conflicts of interest to declare",
    " ",
    "$mojiBullet The study used publicly available de-identified data and did not require IRB approval",
    " ",
    "$mojiBullet All data, code, and analysis scripts will be made publicly available upon acceptance"
)
Replace-Paragraph 21 "BlockText" (Fix-Runs $p21 $fixBulletEq)

# --- Paragraph 22: "COMPETING INTERESTS" heading -> restyled to FirstParagraph (no text change) ---
$d.Paragraphs.Item(22).Style = "FirstParagraph"

# --- Paragraph 30: co-author list with corrupted en dashes (stays BodyText) ---
$p30 = @(
    "Co-authors:",
    " ",
    "Dwi Anggriani $mojiDash Institut Sains Teknologi dan Kesehatan " + [char]0x2019 + "Aisyiyah Kendari",
    " ",
    "Muhammad Atnang $mojiDash Institut Sains Teknologi dan Kesehatan " + [char]0x2019 + "Aisyiyah Kendari",
    " ",
    "Kartini Aprilia Pratiwi Nuzry $mojiDash Institut Sains Teknologi dan Kesehatan " + [char]0x2019 + "Aisyiyah Kendari"
)
Replace-Paragraph 30 "BodyText" (Fix-Runs $p30 $fixBulletPlus)

# --- Paragraph 32: MANUSCRIPT HIGHLIGHTS bullets -> restyled to BlockText ---
$p32 = @(
    "$mojiBullet First systematic evaluation of LLM consistency versus accuracy in medical diagnosis",
    " ",
    "$mojiBullet 1,200 predictions from three state-of-the-art models with rigorous checkpoint system",
    " ",
    "$mojiBullet 99-100% consistency but only 50% accuracy $mojiDash unprecedented 50-point gap",
    " ",
    "$mojiBullet Systematic positive diagnosis bias (49-51 false positives, 0-1 false negatives)",
    " ",
    "$mojiBullet Prompt engineering had minimal effect, suggesting deep-rooted model behavior",
    " ",
    "$mojiBullet Recommends LLMs as supplementary tools, not primary diagnostic systems"
)
Replace-Paragraph 32 "BlockText" (Fix-Runs $p32 $fixBulletEq)

Write-Output "Mojibake fix applied."
